# Weekly update: a new daily-price record is inserted for Choclo at
# "Macroferia Regional de Talca" (Hortaliza sheet). The new record is
# inserted as row 86, and every row that used to be at 86..178 shifts
# down by one (to 87..179). The new row 86 starts out identical to the
# record that used to occupy that slot, then gets its own date (Fecha)
# and volume (Volumen) for the newly reported day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 86 downward (to 87) to make room for the new record.
$ws.Rows(86).Insert()

# Seed the freshly-inserted (now blank) row 86 with the same data the
# row below it holds (that row is the original row 86, now shifted to 87).
$ws.Range("A87:R87").Copy()
$ws.Range("A86").PasteSpecial()

# Overwrite with this week's actual figures: a new date (one day later)
# and the newly reported trading volume.
$ws.Cells.Item(86, 4).Value = 44587   # D86 - Fecha
$ws.Cells.Item(86, 10).Value = 60000  # J86 - Volumen
